$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 794.102088
$schedule.Range("F2").Value = 13.12999484126984

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B12").Value = 57.1
$detailed.Range("B13").Value = 57.3

$detailed.Range("B15").Value = 36.2
$detailed.Range("C15").Value = "historical"

$detailed.Range("C16").Value = "historical"

$detailed.Range("B17").Value = 36.06029
$detailed.Range("B18").Value = 29.57649
$detailed.Range("B19").Value = 36.06011
$detailed.Range("B20").Value = 36.06011

$detailed.Range("B23").Value = 36.06032
$detailed.Range("B24").Value = 36.06046
$detailed.Range("B25").Value = 50.35718
$detailed.Range("B26").Value = 36.06092

$detailed.Range("B30").Value = 36.06
$detailed.Range("B31").Value = 36.06031
$detailed.Range("B32").Value = 30.35616
$detailed.Range("B33").Value = 36.06038
$detailed.Range("B34").Value = 33.26801
$detailed.Range("B35").Value = 8.41405
$detailed.Range("B36").Value = -3.07809
$detailed.Range("B37").Value = -3.01858
$detailed.Range("B38").Value = -2.91738
$detailed.Range("B39").Value = -2.86323
$detailed.Range("B40").Value = 3.4639
$detailed.Range("B41").Value = 9.5329
$detailed.Range("B42").Value = 29.71119
$detailed.Range("B43").Value = 9.59267
$detailed.Range("B44").Value = 8.33337
$detailed.Range("B45").Value = 6.57492
$detailed.Range("B46").Value = 30.05581
$detailed.Range("B47").Value = 56.98
$detailed.Range("B49").Value = 47.61072
